$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Obrigatorio" column (E) for rows 2-6 from "N" to "S"
$ws.Range("E2:E6").Value = "S"
